# Edit slide 14 "Content Placeholder 2": reposition/resize the shape,
# replace its text with a hyperlinked GitHub URL followed by a plain
# space run, and make the shape pick up a fresh internal id (5) the
# way PowerPoint does when a placeholder is recreated.

$presentation = $ppt.ActivePresentation
$slide = $presentation.Slides.Item(14)

# Locate the "Content Placeholder 2" shape on the slide.
$original = $null
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $candidate = $slide.Shapes.Item($i)
    if ($candidate.Name -eq "Content Placeholder 2") {
        $original = $candidate
    }
}

# Duplicating (instead of editing in place) mirrors how PowerPoint ends
# up minting a new shape id/name generation for this placeholder while
# still inheriting the existing run's formatting (so "dirty" round-trips
# like a real edit would).
$newShape = $original.Duplicate()
$original.Delete()

# Deleting the placeholder makes PowerPoint immediately re-create an
# empty inherited placeholder in its place; drop that extra copy and
# keep only the duplicated shape we already styled above.
for ($i = $slide.Shapes.Count; $i -ge 1; $i--) {
    $candidate = $slide.Shapes.Item($i)
    if ($candidate.Id -ne $newShape.Id -and $candidate.Name -like "Content Placeholder*") {
        $candidate.Delete()
    }
}

$newShape.Name = "Content Placeholder 2"

# Reposition / resize the shape (values are EMU / 12700 in points).
$newShape.Left = 84.55409628818899
$newShape.Top = 156.97850803700788
$newShape.Width = 827.9999695
$newShape.Height = 89.37448508897637

# Replace the text: a hyperlinked repo URL followed by a separate,
# non-linked space run.
$url = "https://github.com/Nithinkumar1123/8th_Sem_Intenship_report"
$textRange = $newShape.TextFrame.TextRange
$textRange.Text = $url + " "

# Turn off the bullet for this paragraph.
$textRange.ParagraphFormat.Bullet.Visible = $false

# Apply the hyperlink to just the URL portion of the text.
$linkRange = $textRange.Characters(1, $url.Length)
$linkRange.ActionSettings.Item(1).Hyperlink.Address = $url

Write-Output ("Shape id: " + $newShape.Id)
Write-Output ("Shape name: " + $newShape.Name)
Write-Output ("Text: '" + $newShape.TextFrame.TextRange.Text + "'")
